# Tarea para subir archivos a input type file
# Update the "image" sample value on the Profile sheet to a local-style
# file path, and leave the selection on the edited cell (C2) as Excel
# would after typing a value and confirming it in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Profile")

$ws.Activate()
$ws.Range("C2").Value = "C:/Wappi/Wappi/img.jpg"
$ws.Range("C2").Select() | Out-Null
